$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D24").Value = "아무 글이나 써보고 싶어서 (2)"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/222264870135"

$ws.Range("D32").Value = "한개의 모델로 성격이 비슷한 여러개의 모델을 대체해보자"
$ws.Range("E32").Value = "https://dodonam.tistory.com/299"

$ws.Range("D36").Value = "Applications of Deep Learning for Computer Vision"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/312"

$ws.Range("D39").Value = "Probability concepts explained: Introduction"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Probability-concepts-explained-Introduction-1"

$ws.Range("D51").Value = "[세이버메트릭스] MLB에서 2019까지 1000안타 이상 친 선수 중, 홈런을 2루타보다 더 많이 쳐낸 선수는? (sqlite3)"
$ws.Range("E51").Value = "https://bskyvision.com/1130"

$wb.Save()
